$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$win = $excel.ActiveWindow
Write-Host "ScrollRow=$($win.ScrollRow) ScrollColumn=$($win.ScrollColumn)"
$win.ScrollColumn = 3
$win.ScrollRow = 2
Write-Host "After: ScrollRow=$($win.ScrollRow) ScrollColumn=$($win.ScrollColumn)"
